$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the VRAM memory map boundaries now that it's finished in sheet.c
# (order chosen to match the resulting shared-strings table ordering)
$ws.Range("C4").Value = "9236 B"
$ws.Range("A4").Value = "0x00268000 - 0x0026a413"
$ws.Range("A5").Value = "0x0026a414 - 0x0026f7ff"
$ws.Range("C5").Value = "21484 B"

# Update the sheet view: scroll so row 4 is at top and select A11 instead of D11
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("A11").Select()
